$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 1258.1666
$ws.Range("I32").Value = 1050
$ws.Range("J32").Value = 1299.8
$ws.Range("K32").Value = 1050
$ws.Range("L32").Value = 1299.8
$ws.Range("M32").Value = -724
$ws.Range("N32").Value = -1951.8
$ws.Range("H43").Value = 1565.2084
$ws.Range("J43").Value = 1750.7368
$ws.Range("L43").Value = 1750.7368
$ws.Range("N43").Value = -1888.7368

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 3270.4
$ws.Range("I28").Value = 3270.4
$ws.Range("K28").Value = 3270.4
$ws.Range("M28").Value = -3078.4
$ws.Range("H32").Value = 18615.555
$ws.Range("I32").Value = 17419.451
$ws.Range("J32").Value = 23211.105
$ws.Range("K32").Value = 17419.451
$ws.Range("L32").Value = 23211.105
$ws.Range("M32").Value = -17132.451
$ws.Range("N32").Value = -23785.105
$ws.Range("H61").Value = 52738936
$ws.Range("I61").Value = 71501390
$ws.Range("K61").Value = 71501390
$ws.Range("M61").Value = -71501178
$ws.Range("H99").Value = 3270.4
$ws.Range("I99").Value = 3270.4
$ws.Range("K99").Value = 3270.4
$ws.Range("M99").Value = -275.4000000000001
$ws.Range("H129").Value = 48799.8
$ws.Range("J129").Value = 48799.8
$ws.Range("L129").Value = 48799.8
$ws.Range("N129").Value = -58799.8
$ws.Range("H132").Value = 89843.086
$ws.Range("I132").Value = 52811.8
$ws.Range("J132").Value = 274999.5
$ws.Range("K132").Value = 158435.4
$ws.Range("L132").Value = 824998.5
$ws.Range("M132").Value = -155905.4
$ws.Range("N132").Value = -830058.5
$ws.Range("H136").Value = 52738936
$ws.Range("I136").Value = 71501390
$ws.Range("K136").Value = 214504170
$ws.Range("M136").Value = -214501620

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1300
$ws.Range("I99").Value = 1100
$ws.Range("K99").Value = 1100
$ws.Range("M99").Value = 398
$ws.Range("H122").Value = 1992
$ws.Range("I122").Value = 1748.3636
$ws.Range("J122").Value = 2260
$ws.Range("K122").Value = 5245.0908
$ws.Range("L122").Value = 6780
$ws.Range("M122").Value = -2795.0908
$ws.Range("N122").Value = -11680
$ws.Range("H126").Value = 1300
$ws.Range("I126").Value = 1100
$ws.Range("K126").Value = 3300
$ws.Range("M126").Value = -830

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 571.381
$ws.Range("I5").Value = 240.6923
$ws.Range("J5").Value = 1108.75
$ws.Range("K5").Value = 722.0769
$ws.Range("L5").Value = 3326.25
$ws.Range("M5").Value = -610.0769
$ws.Range("N5").Value = -3550.25
$ws.Range("H113").Value = 634.2759
$ws.Range("I113").Value = 551
$ws.Range("J113").Value = 792.5
$ws.Range("K113").Value = 1653
$ws.Range("L113").Value = 2377.5
$ws.Range("M113").Value = 517
$ws.Range("N113").Value = -6717.5
$ws.Range("H131").Value = 1035.6316
$ws.Range("I131").Value = 408.55554
$ws.Range("J131").Value = 1600
$ws.Range("K131").Value = 1225.66662
$ws.Range("L131").Value = 4800
$ws.Range("M131").Value = 3814.33338
$ws.Range("N131").Value = -14880
$ws.Range("H135").Value = 571.381
$ws.Range("I135").Value = 240.6923
$ws.Range("J135").Value = 1108.75
$ws.Range("K135").Value = 2166.2307
$ws.Range("L135").Value = 9978.75
$ws.Range("M135").Value = 368.7692999999999
$ws.Range("N135").Value = -15048.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 66163.55
$ws.Range("I132").Value = 46865.684
$ws.Range("K132").Value = 140597.052
$ws.Range("M132").Value = -138067.052

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 11908466
$ws.Range("I7").Value = 27780598
$ws.Range("J7").Value = 4366.1665
$ws.Range("K7").Value = 27780598
$ws.Range("L7").Value = 4366.1665
$ws.Range("M7").Value = -27780486
$ws.Range("N7").Value = -4590.1665
$ws.Range("H22").Value = 640.087
$ws.Range("I22").Value = 610.9091
$ws.Range("J22").Value = 666.8333
$ws.Range("K22").Value = 610.9091
$ws.Range("L22").Value = 666.8333
$ws.Range("M22").Value = -315.9091
$ws.Range("N22").Value = -1256.8333
$ws.Range("H27").Value = 640.087
$ws.Range("I27").Value = 610.9091
$ws.Range("J27").Value = 666.8333
$ws.Range("K27").Value = 610.9091
$ws.Range("L27").Value = 666.8333
$ws.Range("M27").Value = -503.9091
$ws.Range("N27").Value = -880.8333
$ws.Range("H40").Value = 2452.5833
$ws.Range("I40").Value = 2754.3333
$ws.Range("J40").Value = 2150.8333
$ws.Range("K40").Value = 2754.3333
$ws.Range("L40").Value = 2150.8333
$ws.Range("M40").Value = -2618.3333
$ws.Range("N40").Value = -2422.8333
$ws.Range("H82").Value = 1166.6666
$ws.Range("I82").Value = 600
$ws.Range("J82").Value = 1450
$ws.Range("K82").Value = 600
$ws.Range("L82").Value = 1450
$ws.Range("M82").Value = -239
$ws.Range("N82").Value = -2172
$ws.Range("H85").Value = 1166.6666
$ws.Range("I85").Value = 600
$ws.Range("J85").Value = 1450
$ws.Range("K85").Value = 600
$ws.Range("L85").Value = 1450
$ws.Range("M85").Value = 648
$ws.Range("N85").Value = -3946
$ws.Range("H99").Value = 15333.333
$ws.Range("I99").Value = 13000
$ws.Range("J99").Value = 20000
$ws.Range("K99").Value = 13000
$ws.Range("L99").Value = 20000
$ws.Range("M99").Value = -10005
$ws.Range("N99").Value = -25990
$ws.Range("H100").Value = 1460.2963
$ws.Range("I100").Value = 1251.75
$ws.Range("J100").Value = 1763.6364
$ws.Range("K100").Value = 1251.75
$ws.Range("L100").Value = 1763.6364
$ws.Range("M100").Value = -710.75
$ws.Range("N100").Value = -2845.6364
$ws.Range("H122").Value = 3358.0728
$ws.Range("I122").Value = 2894.923
$ws.Range("J122").Value = 3501.4285
$ws.Range("K122").Value = 8684.769
$ws.Range("L122").Value = 10504.2855
$ws.Range("M122").Value = -6234.769
$ws.Range("N122").Value = -15404.2855
$ws.Range("H123").Value = 24571.428
$ws.Range("J123").Value = 24571.428
$ws.Range("L123").Value = 24571.428
$ws.Range("N123").Value = -34371.428
$ws.Range("H125").Value = 45000
$ws.Range("J125").Value = 45000
$ws.Range("L125").Value = 45000
$ws.Range("N125").Value = -54840
$ws.Range("H126").Value = 11908466
$ws.Range("I126").Value = 27780598
$ws.Range("J126").Value = 4366.1665
$ws.Range("K126").Value = 83341794
$ws.Range("L126").Value = 13098.4995
$ws.Range("M126").Value = -83339324
$ws.Range("N126").Value = -18038.4995
$ws.Range("H127").Value = 49508.25
$ws.Range("J127").Value = 49508.25
$ws.Range("L127").Value = 49508.25
$ws.Range("N127").Value = -59428.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1607.5385
$ws.Range("I126").Value = 1694.4
$ws.Range("J126").Value = 1318
$ws.Range("K126").Value = 5083.200000000001
$ws.Range("L126").Value = 3954
$ws.Range("M126").Value = -2613.200000000001
$ws.Range("N126").Value = -8894
$ws.Range("H136").Value = 58924.715
$ws.Range("I136").Value = 38239.445
$ws.Range("J136").Value = 128737.5
$ws.Range("K136").Value = 114718.335
$ws.Range("L136").Value = 386212.5
$ws.Range("M136").Value = -112168.335
$ws.Range("N136").Value = -391312.5
